# Mise à jour de l'application
# Adds a new training-day column (CP, 16/12/2025) after the existing last
# column (CO, 12/12/2025) on the attendance sheet, mirroring the style of
# the previous date column and filling in each player's attendance mark.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New date header in CP1, same formatting as the previous date cell (CO1).
$ws.Cells.Item(1, 94).Value = 46007
$ws.Cells.Item(1, 93).Copy()
$ws.Cells.Item(1, 94).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Attendance marks for the new date, one per player row. Rows absent from
# this table are left untouched (row 12's tracking stopped earlier in the
# sheet, so it gets no new cell); row 21 gets a styled-but-empty cell since
# that player already has no attendance data recorded.
$marks = @{
    2  = "P"
    3  = "P"
    4  = "P"
    5  = "P"
    6  = "P"
    7  = "P"
    8  = "P"
    9  = "P"
    10 = "P"
    11 = "P"
    13 = "B"
    14 = "P"
    15 = "P"
    16 = "B"
    17 = "P"
    18 = "P"
    19 = "P"
    20 = "P"
    22 = "P"
    23 = "B"
    24 = "P"
    25 = "P"
    26 = "P"
    27 = "P"
    28 = "P"
    29 = "P"
}

foreach ($row in $marks.Keys) {
    $ws.Cells.Item($row, 94).Value = $marks[$row]
    $ws.Cells.Item($row, 93).Copy()
    $ws.Cells.Item($row, 94).PasteSpecial(-4122)
    $excel.CutCopyMode = $false
}

# Row 21: styled empty cell matching the rest of that (dataless) row.
$ws.Cells.Item(21, 93).Copy()
$ws.Cells.Item(21, 94).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Leave the cursor where the author's session ended up.
$ws.Range("CR24").Select() | Out-Null
